$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '24.315.25'
Set-TextValue $ws.Range('E2') '  -5.72%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.630.60'
Set-TextValue $ws.Range('E3') '  -7.34%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.007'
Set-TextValue $ws.Range('E4') '  +0.23%  '

# Row 5
Set-TextValue $ws.Range('D5') '1.003'
Set-TextValue $ws.Range('E5') '  +0.29%  '

# Row 6
Set-TextValue $ws.Range('D6') '304.15'

# Row 7
Set-TextValue $ws.Range('D7') '0.3606'
Set-TextValue $ws.Range('E7') '  -5.52%  '

# Row 8
Set-TextValue $ws.Range('D8') '46.81'
Set-TextValue $ws.Range('E8') '  -6.95%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.3218'
Set-TextValue $ws.Range('E9') '  -10.59%  '

# Row 10
Set-TextValue $ws.Range('D10') '1.097'
Set-TextValue $ws.Range('E10') '  -10.41%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.06853'
Set-TextValue $ws.Range('E11') '  -10.72%  '

# Row 12
Set-TextValue $ws.Range('D12') '1.005'
Set-TextValue $ws.Range('E12') '  +0.38%  '

# Row 13
Set-TextValue $ws.Range('D13') '5.893'
Set-TextValue $ws.Range('E13') '  -8.86%  '

# Row 14
Set-TextValue $ws.Range('D14') '19.01'
Set-TextValue $ws.Range('E14') '  -12.10%  '

# Row 15
Set-TextValue $ws.Range('D15') '1.636.04'
Set-TextValue $ws.Range('E15') '  -7.25%  '

# Row 16
Set-TextValue $ws.Range('D16') '6.500'
Set-TextValue $ws.Range('E16') '  -8.20%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.00001038'
Set-TextValue $ws.Range('E17') '  -9.96%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.06512'
Set-TextValue $ws.Range('E18') '  -4.08%  '

# Row 19
Set-TextValue $ws.Range('D19') '1.003'
Set-TextValue $ws.Range('E19') '  +0.34%  '

# Row 20
Set-TextValue $ws.Range('D20') '76.08'
Set-TextValue $ws.Range('E20') '  -12.36%  '

# Row 21
Set-TextValue $ws.Range('B21') 'Uniswap'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D21') '5.855'
Set-TextValue $ws.Range('E21') '  -9.97%  '

# Row 22
Set-TextValue $ws.Range('B22') 'Avalanche'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D22') '15.63'
Set-TextValue $ws.Range('E22') '  -11.45%  '

# Row 23
Set-TextValue $ws.Range('D23') '11.92'
Set-TextValue $ws.Range('E23') '  -8.35%  '

# Row 24
Set-TextValue $ws.Range('D24') '24.310.59'
Set-TextValue $ws.Range('E24') '  -5.71%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.391'

# Row 26
Set-TextValue $ws.Range('D26') '2.299'
Set-TextValue $ws.Range('E26') '  -20.67%  '

# Row 27
Set-TextValue $ws.Range('D27') '143.67'
Set-TextValue $ws.Range('E27') '  -7.83%  '

# Row 28
Set-TextValue $ws.Range('D28') '18.46'
Set-TextValue $ws.Range('E28') '  -11.04%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.818.80'
Set-TextValue $ws.Range('E29') '  -7.21%  '

# Row 30
Set-TextValue $ws.Range('D30') '123.56'
Set-TextValue $ws.Range('E30') '  -7.62%  '

# Row 31
Set-TextValue $ws.Range('D31') '1.105'
Set-TextValue $ws.Range('E31') '  -8.37%  '

# Row 32
Set-TextValue $ws.Range('D32') '4.073'
Set-TextValue $ws.Range('E32') '  -3.31%  '

# Row 33
Set-TextValue $ws.Range('D33') '5.590'
Set-TextValue $ws.Range('E33') '  -22.15%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.08344'
Set-TextValue $ws.Range('E34') '  -4.77%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.663'
Set-TextValue $ws.Range('E35') '  -8.02%  '

# Row 36
Set-TextValue $ws.Range('D36') '12.21'
Set-TextValue $ws.Range('E36') '  -14.59%  '

# Row 37
Set-TextValue $ws.Range('D37') '5.075'
Set-TextValue $ws.Range('E37') '  -11.31%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.05968'
Set-TextValue $ws.Range('E38') '  -11.59%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.02198'
Set-TextValue $ws.Range('E39') '  -11.98%  '

# Row 40
Set-TextValue $ws.Range('D40') '1.198'
Set-TextValue $ws.Range('E40') '  -7.55%  '

# Row 41
Set-TextValue $ws.Range('D41') '8.090'
Set-TextValue $ws.Range('E41') '  -13.70%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.2019'
Set-TextValue $ws.Range('E42') '  -10.64%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.003'
Set-TextValue $ws.Range('E43') '  +0.39%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.5798'
Set-TextValue $ws.Range('E44') '  -11.92%  '

# Row 45
Set-TextValue $ws.Range('D45') '3.702'
Set-TextValue $ws.Range('E45') '  -5.15%  '

# Row 46
Set-TextValue $ws.Range('D46') '12.47'
Set-TextValue $ws.Range('E46') '  -13.65%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.5485'
Set-TextValue $ws.Range('E47') '  -13.48%  '

# Row 48
Set-TextValue $ws.Range('D48') '120.92'
Set-TextValue $ws.Range('E48') '  -8.46%  '

# Row 49
Set-TextValue $ws.Range('D49') '1.906'

# Row 50
Set-TextValue $ws.Range('D50') '0.06852'
Set-TextValue $ws.Range('E50') '  -8.66%  '

# Row 51
Set-TextValue $ws.Range('D51') '72.92'
Set-TextValue $ws.Range('E51') '  -9.91%  '
